# Fixed variables and query errors in Bread from TC01 to TC30
#
# The "Cases" query stored in cell B2 of the "startup" sheet contained a
# trailing line that referenced an undeclared/unused `co.cohort_description`
# (Cohort) column. That line (and its trailing blank line) is removed here,
# matching the corrected Cypher query used across the TC01-TC30 test files.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Bernese Mountain Dog']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

# Normalize to a single trailing-newline-free string (here-string captures one
# extra newline before the closing '@ marker).
$newCasesQuery = $newCasesQuery.TrimEnd("`r", "`n")

$ws.Range("B2").Value = $newCasesQuery

# Refresh the view: the workbook was reopened/scrolled back to the top of the
# sheet and re-zoomed before saving.
$aw = $excel.ActiveWindow
$aw.Zoom = 115
$ws.Range("B2").Select()
